$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$vals = @(1.834436176833037,1.787766169578248,1.760161441781094,1.749176177257425,1.747368010421667,1.760012222666631,1.818126017185847,1.94044999472203,2.035467219018017,2.079821624293686,2.096780739631924,2.09312103043834,2.081213590015523,2.073941189489119,2.032591345993353,2.00751446969025,1.993197300436691,1.988368010375495,2.010172931901025,2.08470666442139,2.134369470329602,2.10777634176452,2.008970730118904,1.906457948726597)
$arr = New-Object 'object[,]' $vals.Length,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("B2:B25").Value = $arr

$vals = @(0.02932362838556912,0.02551690567669596,0.02317775130687494,0.02222398558131999,0.02206557831810585,0.02316489079306905,0.02801139342757608,0.03750531618781849,0.04448098776195764,0.04765597956030376,0.04885862066682023,0.04859959375117739,0.04775491420666356,0.04723757072378021,0.0442735371751013,0.04245570579429625,0.04141029391666962,0.04105635951142972,0.0426492001654708,0.04800300707665883,0.05150405849049378,0.04963526490865888,0.04256172235393763,0.03493735427075251)
$arr = New-Object 'object[,]' $vals.Length,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("C2:C25").Value = $arr

$vals = @(0.06239476398544497,0.06185829364868312,0.06155810626455249,0.06144312583019662,0.06142447735992285,0.06155652583936799,0.06220372892349957,0.06370462145329014,0.06494879307343382,0.06554557804981087,0.06577599586443839,0.06572617433576866,0.06556444591966581,0.06546595925197352,0.0649104115197261,0.06457748979104494,0.06438890160638522,0.06432554694268333,0.06461262979202687,0.06561182930331455,0.06629067496988839,0.06592600095524404,0.06459673423584178,0.06327376259238804)
$arr = New-Object 'object[,]' $vals.Length,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("E2:E25").Value = $arr

$vals = @(0.4443680307746263,0.387822817061874,0.3531389305168915,0.3390132514313251,0.336668177824194,0.3529483938344953,0.4248636149813336,0.5661985755041457,0.6702781546542269,0.7176906081379002,0.7356546913071611,0.7317853510981394,0.7191683204515869,0.7114413442032514,0.6671810134426437,0.6400460337215605,0.6244449056556647,0.6191636801734006,0.642933953830422,0.7228739723492197,0.7751780083420101,0.7472568307915566,0.6416283278902313,0.5279251897347308)
$arr = New-Object 'object[,]' $vals.Length,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("F2:F25").Value = $arr

$vals = @(0.002607791020105874,0.00261235509399859,0.002615304096617995,0.002616542835987643,0.00261675076612204,0.002615320652651072,0.00260933435114246,0.002598753171530571,0.002591677226515037,0.002588608091820887,0.002587467296375287,0.002587712036496815,0.002588513809271143,0.002589007704425223,0.002591880804972818,0.002593681631648132,0.002594731521386214,0.002595089421598125,0.002593488471824081,0.002588277728806924,0.002584997001358508,0.002586736606439944,0.002593575754365893,0.002601492510990775)
$arr = New-Object 'object[,]' $vals.Length,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("G2:G25").Value = $arr

$vals = @(1.324321461773508,1.27890035438287,1.251807035296622,1.240965817470453,1.239177678687867,1.251660019592435,1.308494860281456,1.426288999078452,1.516750877236404,1.558767459013723,1.574803215903302,1.571344060094731,1.560084221922779,1.553203541277071,1.514022452215642,1.490208098179608,1.476592094585357,1.471995909260215,1.492734749696439,1.56338811230674,1.610292991045725,1.585192085103216,1.491592216399312,1.393738291804794)
$arr = New-Object 'object[,]' $vals.Length,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("K2:K25").Value = $arr

$vals = @(0.2466298486510681,0.2399281671627875,0.2359576683570168,0.234375893155871,0.2341154272456833,0.2359361893522305,0.2442891042623074,0.2618187712697164,0.275406250624485,0.2817431503441412,0.2841652964216337,0.2836426418982541,0.281941970347674,0.2809031925480951,0.2749952632464669,0.271410925035056,0.269363976233592,0.2686734312347312,0.2717909653034383,0.2824408874984385,0.2895324128366212,0.2857355007445079,0.2716191063238682,0.2569526793201646)
$arr = New-Object 'object[,]' $vals.Length,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("L2:L25").Value = $arr

$vals = @(3.810167952200473,3.809469029033608,3.809768174239593,3.810072772250493,3.810134371993826,3.809771543043908,3.809775538144251,3.815583654044758,3.823419489406888,3.827766386871303,3.829525500019855,3.829141607510195,3.827908841625799,3.82716847381738,3.823151197921845,3.820887485483297,3.819659062640113,3.81925576512981,3.821120839905888,3.828267862721844,3.8335979782903,3.830692705934041,3.821015112957042,3.81338809684172)
$arr = New-Object 'object[,]' $vals.Length,1
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[$i,0] = $vals[$i] }
$ws.Range("N2:N25").Value = $arr

